# sku_wise_target_sales.xlsx: rename headers, refresh target/sales figures,
# and re-sort SKU rows to match the new Item-code ordering.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header relabeling ---
$ws.Range("A1").Value = "Item"
$ws.Range("B1").Value = "skuname"
$ws.Range("C1").Value = "Months Target(Tk)"
$ws.Range("D1").Value = "MTD Target(Tk)"
$ws.Range("E1").Value = "MTD Sales(Tk)"
$ws.Range("F1").Value = "Months Target(Kg)"
$ws.Range("G1").Value = "MTD Target(Kg)"
$ws.Range("H1").Value = "MTD Sales(Kg)"
$ws.Range("I1").Value = "Achiv%(Tk)"
$ws.Range("J1").Value = "Achiv%(Kg)"

# --- Rows 2-27: SKU data (re-ordered by Item code, values refreshed) ---
$data = New-Object 'object[,]' 26,10
# row 2: Kurkure (20 gm) - Chilli Chatka
$data[0,0] = 142
$data[0,1] = "Kurkure (20 gm) - Chilli Chatka"
$data[0,2] = 20866.65
$data[0,3] = 13462
$data[0,4] = 4287.6
$data[0,5] = 50.1
$data[0,6] = 32
$data[0,7] = 10.8
$data[0,8] = 31.84965086911306
$data[0,9] = 33.75
# row 3: Kurkure (20 gm) - Masala Munch
$data[1,0] = 143
$data[1,1] = "Kurkure (20 gm) - Masala Munch"
$data[1,2] = 20866.65
$data[1,3] = 13462
$data[1,4] = 1381.56
$data[1,5] = 50.1
$data[1,6] = 32
$data[1,7] = 3.48
$data[1,8] = 10.26266528004754
$data[1,9] = 10.875
# row 4: Kurkure (20 gm) - Tock Misti Jhal
$data[2,0] = 144
$data[2,1] = "Kurkure (20 gm) - Tock Misti Jhal"
$data[2,2] = 12519.99
$data[2,3] = 8077
$data[2,4] = 666.96
$data[2,5] = 30.06
$data[2,6] = 19
$data[2,7] = 1.68
$data[2,8] = 8.257521356939458
$data[2,9] = 8.842105263157894
# row 5: Kurkure (20 gm)- ASCO
$data[3,0] = 146
$data[3,1] = "Kurkure (20 gm)- ASCO"
$data[3,2] = 62608.28
$data[3,3] = 40392
$data[3,4] = 5335.68
$data[3,5] = 150.32
$data[3,6] = 97
$data[3,7] = 13.44
$data[3,8] = 13.20974450386215
$data[3,9] = 13.85567010309278
# row 6: Kurkure (90 gm) - Chilli Chatka
$data[4,0] = 154
$data[4,1] = "Kurkure (90 gm) - Chilli Chatka"
$data[4,2] = 4208.67
$data[4,3] = 2715
$data[4,4] = 1944.32
$data[4,5] = 9.09
$data[4,6] = 6
$data[4,7] = 4.41
$data[4,8] = 71.61399631675874
$data[4,9] = 73.5
# row 7: Kurkure (90 gm) - Masala Munch
$data[5,0] = 155
$data[5,1] = "Kurkure (90 gm) - Masala Munch"
$data[5,2] = 4208.67
$data[5,3] = 2715
$data[5,4] = 2301.44
$data[5,5] = 9.09
$data[5,6] = 6
$data[5,7] = 5.22
$data[5,8] = 84.76758747697974
$data[5,9] = 87
# row 8: Kurkure (90 gm) - Tock Misti Jhal
$data[6,0] = 156
$data[6,1] = "Kurkure (90 gm) - Tock Misti Jhal"
$data[6,2] = 2541.87
$data[6,3] = 1640
$data[6,4] = 79.36
$data[6,5] = 5.49
$data[6,6] = 4
$data[6,7] = 0.18
$data[6,8] = 4.839024390243902
$data[6,9] = 4.5
# row 9: Kurkure (90 gm) - ASCO
$data[7,0] = 157
$data[7,1] = "Kurkure (90 gm) - ASCO"
$data[7,2] = 12667.68
$data[7,3] = 8173
$data[7,4] = 3690.24
$data[7,5] = 27.36
$data[7,6] = 18
$data[7,7] = 8.37
$data[7,8] = 45.15159672091031
$data[7,9] = 46.5
# row 10: Kurkure (20 gm) - STT
$data[8,0] = 158
$data[8,1] = "Kurkure (20 gm) - STT"
$data[8,2] = 50088.29
$data[8,3] = 32315
$data[8,4] = 2096.16
$data[8,5] = 120.26
$data[8,6] = 78
$data[8,7] = 5.28
$data[8,8] = 6.486647067925112
$data[8,9] = 6.769230769230769
# row 11: Kurkure (45 gm) - STT
$data[9,0] = 159
$data[9,1] = "Kurkure (45 gm) - STT"
$data[9,2] = 23892.01
$data[9,3] = 15414
$data[9,4] = 3809.28
$data[9,5] = 51.615
$data[9,6] = 33
$data[9,7] = 8.64
$data[9,8] = 24.71311794472558
$data[9,9] = 26.18181818181818
# row 12: Kurkure (90 gm) - STT
$data[10,0] = 160
$data[10,1] = "Kurkure (90 gm) - STT"
$data[10,2] = 10167.48
$data[10,3] = 6560
$data[10,4] = 3253.76
$data[10,5] = 21.96
$data[10,6] = 14
$data[10,7] = 7.38
$data[10,8] = 49.60000000000001
$data[10,9] = 52.71428571428572
# row 13: Lays(25 gm) ASCO
$data[11,0] = 166
$data[11,1] = "Lays(25 gm) ASCO"
$data[11,2] = 8221.5
$data[11,3] = 5304
$data[11,4] = 6582.6
$data[11,5] = 9.45
$data[11,6] = 6
$data[11,7] = 7.95
$data[11,8] = 124.106334841629
$data[11,9] = 132.5
# row 14: Lays (25 gm) - STT
$data[12,0] = 168
$data[12,1] = "Lays (25 gm) - STT"
$data[12,2] = 4741.5
$data[12,3] = 3059
$data[12,4] = 7493.4
$data[12,5] = 5.45
$data[12,6] = 4
$data[12,7] = 9.05
$data[12,8] = 244.9624060150376
$data[12,9] = 226.25
# row 15: Lays(52 gm) ASCO
$data[13,0] = 170
$data[13,1] = "Lays(52 gm) ASCO"
$data[13,2] = 167692.5
$data[13,3] = 108189
$data[13,4] = 64019.86
$data[13,5] = 200.46
$data[13,6] = 129
$data[13,7] = 80.392
$data[13,8] = 59.17409348455018
$data[13,9] = 62.31937984496124
# row 16: Lays (52 gm) - STT
$data[14,0] = 172
$data[14,1] = "Lays (52 gm) - STT"
$data[14,2] = 96744
$data[14,3] = 62415
$data[14,4] = 117521.58
$data[14,5] = 115.648
$data[14,6] = 75
$data[14,7] = 147.576
$data[14,8] = 188.2906032203797
$data[14,9] = 196.768
# row 17: Quaker Oats (1000gm)
$data[15,0] = 187
$data[15,1] = "Quaker Oats (1000gm)"
$data[15,2] = 0
$data[15,3] = 0
$data[15,4] = 132425.44
$data[15,5] = 0
$data[15,6] = 0
$data[15,7] = 266
$data[15,8] = 0
$data[15,9] = 0
# row 18: Kurkure Large (45 gm) - ASCO
$data[16,0] = 195
$data[16,1] = "Kurkure Large (45 gm) - ASCO"
$data[16,2] = 29870.22
$data[16,3] = 19271
$data[16,4] = 13332.48
$data[16,5] = 64.53
$data[16,6] = 42
$data[16,7] = 30.24
$data[16,8] = 69.18416273156555
$data[16,9] = 72
# row 19: Kurkure Large (45 gm) - CC
$data[17,0] = 196
$data[17,1] = "Kurkure Large (45 gm) - CC"
$data[17,2] = 9956.74
$data[17,3] = 6424
$data[17,4] = 4920.32
$data[17,5] = 21.51
$data[17,6] = 14
$data[17,7] = 11.16
$data[17,8] = 76.59277708592776
$data[17,9] = 79.71428571428572
# row 20: Kurkure Large (45 gm) - MM
$data[18,0] = 197
$data[18,1] = "Kurkure Large (45 gm) - MM"
$data[18,2] = 9956.74
$data[18,3] = 6424
$data[18,4] = 9602.56
$data[18,5] = 21.51
$data[18,6] = 14
$data[18,7] = 21.78
$data[18,8] = 149.4794520547945
$data[18,9] = 155.5714285714286
# row 21: Kurkure Large (45 gm) - TJM
$data[19,0] = 198
$data[19,1] = "Kurkure Large (45 gm) - TJM"
$data[19,2] = 5978.21
$data[19,3] = 3857
$data[19,4] = 1507.84
$data[19,5] = 12.915
$data[19,6] = 8
$data[19,7] = 3.42
$data[19,8] = 39.0935960591133
$data[19,9] = 42.75
# row 22: Quaker Oats (500 gm Jar)
$data[20,0] = 199
$data[20,1] = "Quaker Oats (500 gm Jar)"
$data[20,2] = 56840.35
$data[20,3] = 36671
$data[20,4] = 8186.17
$data[20,5] = 102.5
$data[20,6] = 66
$data[20,7] = 15.5
$data[20,8] = 22.32327997600284
$data[20,9] = 23.48484848484848
# row 23: Lays 3D-37 gm (Poly)
$data[21,0] = 213
$data[21,1] = "Lays 3D-37 gm (Poly)"
$data[21,2] = 24003.84
$data[21,3] = 15486
$data[21,4] = 5613.3
$data[21,5] = 49.728
$data[21,6] = 32
$data[21,7] = 12.21
$data[21,8] = 36.24757845796204
$data[21,9] = 38.15625
# row 24: Lays Pastazz - 37 gm (Poly)
$data[22,0] = 215
$data[22,1] = "Lays Pastazz - 37 gm (Poly)"
$data[22,2] = 48025.54
$data[22,3] = 30984
$data[22,4] = 4043.2031
$data[22,5] = 99.493
$data[22,6] = 64
$data[22,7] = 8.732
$data[22,8] = 13.04932578104829
$data[22,9] = 13.64375
# row 25: Lays Pastazz - 20 gm
$data[23,0] = 216
$data[23,1] = "Lays Pastazz - 20 gm"
$data[23,2] = 42104.95
$data[23,3] = 27164
$data[23,4] = 4250
$data[23,5] = 84.87
$data[23,6] = 55
$data[23,7] = 9
$data[23,8] = 15.64570755411574
$data[23,9] = 16.36363636363636
# row 26: Lays 3D - 20 gm
$data[24,0] = 217
$data[24,1] = "Lays 3D - 20 gm"
$data[24,2] = 21056.94
$data[24,3] = 13585
$data[24,4] = 510
$data[24,5] = 42.444
$data[24,6] = 27
$data[24,7] = 1.08
$data[24,8] = 3.754140596245859
$data[24,9] = 4
# row 27: Quaker Oats 1Kg New
$data[25,0] = 218
$data[25,1] = "Quaker Oats 1Kg New"
$data[25,2] = 386297.47
$data[25,3] = 249224
$data[25,4] = 11948.184
$data[25,5] = 739
$data[25,6] = 477
$data[25,7] = 24
$data[25,8] = 4.794154656052386
$data[25,9] = 5.031446540880504

$ws.Range("A2:J27").Value = $data
